# river update May 2024
# Append 3 new sample rows (17-19) for "Retaruke at Whanganui Confluence",
# dated 2023-04-06 (serial 45022), mirroring the layout of the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$siteName = "Retaruke at Whanganui Confluence"
$dateSerial = 45022
$dateFormat = $ws.Cells.Item(2, 3).NumberFormat   # same format as the other "date time" cells

$newRows = @(
    @{ Row = 17; Param = "ASPM (Macroinvertebrate Average Score Per Metric)"; Value = "0.532" },
    @{ Row = 18; Param = "MCI (Macroinvertebrate Community Index)";          Value = "112.38" },
    @{ Row = 19; Param = "QMCI (Quantitative Macroinvertebrate Community Index)"; Value = "6.418" }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $siteName
    $ws.Cells.Item($row, 2).Value = $r.Param

    $ws.Cells.Item($row, 3).Value = $dateSerial
    $ws.Cells.Item($row, 3).NumberFormat = $dateFormat

    # Force the numeric-looking value into text (like the source data),
    # then reset the style so no extra number format lingers on the cell.
    $ws.Cells.Item($row, 4).NumberFormat = "@"
    $ws.Cells.Item($row, 4).Value = $r.Value
    $ws.Cells.Item($row, 4).Style = "Normal"

    # Project / Method / Unit / pH columns are blank for every sample row.
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 6).Value = ""
    $ws.Cells.Item($row, 7).Value = ""

    $ws.Cells.Item($row, 8).Value = 200

    $ws.Cells.Item($row, 9).Value = ""
}
